$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price updates in column D (Price), preserving the original text/string cell type
# (cells are stored as text, e.g. "245.08", not numbers) by forcing Text number format
# before assigning, which mirrors how the source feed writes these inline strings.
$priceUpdates = @{
    2 = "245.12"
    3 = "25.53"
    4 = "5.102"
    5 = "0.05587"
    9 = "0.8407"
    10 = "0.1338"
    11 = "0.02866"
    12 = "0.09377"
    14 = "0.0005961"
    15 = "0.006206"
    17 = "2.498"
    19 = "0.06947"
    20 = "0.03194"
    22 = "3.747"
    23 = "0.04703"
    25 = "0.001246"
    26 = "0.004620"
    41 = "0.1351"
    42 = "0.002518"
    43 = "0.003375"
    44 = "0.007511"
    45 = "0.00005321"
    48 = "0.002122"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

